$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Neetcode 150")

# Add the hyperlink first - this stamps the `display` attribute (= URL) on the
# <hyperlink> element. It also temporarily writes the display text into the
# cell and applies the built-in "Hyperlink" style, both of which we overwrite
# below to match the row-24 pattern (plain text value + Neutral/wrap style).
$ws.Hyperlinks.Add($ws.Range("C25"), "https://leetcode.com/problems/letter-combinations-of-a-phone-number/", "", "", "https://leetcode.com/problems/letter-combinations-of-a-phone-number/")

# New row 25 data: Backtracking / Medium / 17. Letter Combinations of a Phone Number
$ws.Range("A25").Value = "Backtracking"
$ws.Range("B25").Value = "Medium"
$ws.Range("C25").Value = "17. Letter Combinations of a Phone Number"
$ws.Range("D25").Value = "Store the letters in an array like [""abc"", ""def"", ...] and convert digits (eg. ""23"" to 01) so that we can correspond it with the array easier`nRun a recursive dfs, at each level, pick each digit in a loop and recurse to next level, until we reach base case i >= n, then append to res with """".join(curr) and then pop out the digit after using it"

# Match the formatting of the row above it (row 24)
$ws.Range("B25").Style = $ws.Range("B24").Style
$ws.Range("C25").Style = $ws.Range("C24").Style
$ws.Range("D25").WrapText = $true
$ws.Range("D25").VerticalAlignment = -4160

$ws.Rows.Item(25).RowHeight = 43.2

$ws.Range("D26").Select()
